$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.059.79'
$ws.Range("E2").Value = '  -1.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.179.23'
$ws.Range("E3").Value = '  -1.37%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.17'
$ws.Range("E5").Value = '  +2.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.620'
$ws.Range("E6").Value = '  -1.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '67.35'
$ws.Range("E7").Value = '  -4.48%  '

$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.570'
$ws.Range("E9").Value = '  +4.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.09'
$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.30'
$ws.Range("E11").Value = '  +0.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0927'
$ws.Range("E12").Value = '  -2.02%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  -0.95%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.98'
$ws.Range("E14").Value = '  +4.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.503.95'
$ws.Range("E15").Value = '  -1.21%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.39'
$ws.Range("E16").Value = '  -2.46%  '

$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.859'
$ws.Range("E17").Value = '  +2.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.172.53'
$ws.Range("E18").Value = '  -1.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '40.957.18'
$ws.Range("E19").Value = '  -1.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0942'
$ws.Range("E20").Value = '  -1.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.18'
$ws.Range("E21").Value = '  +1.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.55'
$ws.Range("E22").Value = '  -2.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.96'
$ws.Range("E23").Value = '  -1.28%  '

$ws.Range("E24").Value = '  +1.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.72'
$ws.Range("E25").Value = '  +19.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.80'
$ws.Range("E26").Value = '  +5.65%  '

$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.47'
$ws.Range("E28").Value = '  +1.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.73'
$ws.Range("E29").Value = '  -3.61%  '

$ws.Range("E30").Value = '  -1.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.77'
$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.48'
$ws.Range("E32").Value = '  +0.37%  '

$ws.Range("E33").Value = '  -1.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.123'
$ws.Range("E34").Value = '  -1.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.39'
$ws.Range("E35").Value = '  +6.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0721'
$ws.Range("E36").Value = '  +1.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.55'
$ws.Range("E37").Value = '  -1.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.36'
$ws.Range("E38").Value = '  +8.87%  '

$ws.Range("E39").Value = '  +2.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0295'
$ws.Range("E40").Value = '  +8.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.21'
$ws.Range("E41").Value = '  -2.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.04'
$ws.Range("E42").Value = '  +17.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.63'
$ws.Range("E43").Value = '  -3.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.16'
$ws.Range("E44").Value = '  -1.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.200'
$ws.Range("E45").Value = '  +4.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.79'
$ws.Range("E46").Value = '  -0.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.101'
$ws.Range("E47").Value = '  +2.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.55'
$ws.Range("E48").Value = '  -4.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.13'
$ws.Range("E50").Value = '  +3.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.15'
$ws.Range("E51").Value = '  -1.48%  '
